$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "staging" URL from A4 up to A2, and drop the long
# Lorem-ipsum placeholder text that used to live in A5.
$ws.Range("A2").Value = $ws.Range("A4").Value()
$ws.Range("A4").ClearContents()
$ws.Range("A5").ClearContents()

# Update the selected cell to match the saved view state.
$ws.Range("A13").Select()
